$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: AHC30716 / CERT III IN HORTICULTURE ---
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("D2").Value = "CERTIFICATE III IN HORTICULTURE"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I2").Value = 10200
$ws.Range("J3").Value = "12,000 tuition fee + 200 handling fee"
$ws.Range("J2").Value = "10,000 tuition fee + 200 handling fee"
$ws.Range("M2").Value = "TAS"

# --- Row 3: AHC40416 / CERT IV IN HORTICULTURE ---
$ws.Range("A3").Value = "AHC40416"
$ws.Range("B3").Value = "110598E"
$ws.Range("D3").Value = "CERTIFICATE IV IN HORTICULTURE"
$ws.Range("E3").Value = 52
$ws.Range("H3").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I3").Value = 12200
$ws.Range("M3").Value = "TAS"

# --- Row 4: AHC51422 / DIPLOMA OF AGRIBUSINESS MANAGEMENT ---
$ws.Range("A4").Value = "AHC51422"
$ws.Range("B4").Value = "110774E"
$ws.Range("D4").Value = "DIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E4").Value = 52
$ws.Range("H4").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("I4").Value = 13200
$ws.Range("J4").Value = "13,000 tuition fee + 200 handling fee"
$ws.Range("M4").Value = "TAS"

# --- Row 5: package CERT III + CERT IV ---
$ws.Range("A5").Value = "AHC30716 / AHC40416"
$ws.Range("B5").Value = "110597F / 110598E"
$ws.Range("D5").Value = "CERTIFICATE III IN HORTICULTURE +`nCERTIFICATE IV IN HORTICULTURE"
$ws.Range("E5").Value = 104
$ws.Range("H5").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I5").Value = 20200
$ws.Range("J5").Value = "20,000 tuition fee + 200 handling fee"
$ws.Range("M5").Value = "TAS"

# --- Row 6: package CERT III + DIPLOMA ---
$ws.Range("A6").Value = "AHC30716 / AHC51422"
$ws.Range("B6").Value = "110597F / 110774E"
$ws.Range("D6").Value = "CERTIFICATE III IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E6").Value = 104
$ws.Range("H6").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I6").Value = 21200
$ws.Range("J6").Value = "21,000 tuition fee + 200 handling fee"
$ws.Range("M6").Value = "TAS"

# --- Row 7: package CERT IV + DIPLOMA ---
$ws.Range("A7").Value = "AHC40416 / AHC51422"
$ws.Range("B7").Value = "110598E / 110774E"
$ws.Range("D7").Value = "CERTIFICATE IV IN HORTICULTURE +`nDIPLOMA OF AGRIBUSINESS MANAGEMENT"
$ws.Range("E7").Value = 104
$ws.Range("H7").Value = "88 wks Tuition + 16 wks Break"
$ws.Range("I7").Value = 21200
$ws.Range("J7").Value = "21,000 tuition fee + 200 handling fee"
$ws.Range("M7").Value = "TAS"

# --- department column, filled last ---
$ws.Range("C2").Value = "HORTICULTURE"
$ws.Range("C3").Value = "HORTICULTURE"
$ws.Range("C4").Value = "MANAGEMENT"
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("C6").Value = "PACKAGES"
$ws.Range("C7").Value = "PACKAGES"

# --- row heights ---
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(7).RowHeight = 45

# --- styles: durationDetail wraps, tuition is formatted as number, tuitionDetail wraps + formatted ---
for ($r = 2; $r -le 7; $r++) {
    $ws.Range("H$r").WrapText = $true
    $ws.Range("I$r").NumberFormat = "#,##0"
    $ws.Range("J$r").NumberFormat = "#,##0"
    $ws.Range("J$r").WrapText = $true
}

# package rows also wrap the multi-code / multi-name columns
$ws.Range("A5:B7").WrapText = $true
$ws.Range("D5:D7").WrapText = $true

$excel.Goto($ws.Range("C10"))
